$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell E8 text from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active selection on the sheet (Excel records the last
# selected cell in the view when the file is saved)
$ws.Activate()
$ws.Range("E8").Select()
